# Add a new "2021" column (column R) to the table, mirroring the
# existing "2020" column (Q) for styling, then fill in the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the formatting of the existing last column (Q3:Q33) into the
#    new column (R3:R33) so borders/number-formats/fonts line up with
#    the rest of the table.
$ws.Range("Q3:Q33").Copy()
$ws.Range("R3:R33").PasteSpecial(-4122)

# 2) Write the new 2021 values (column R) row by row.
$newValues = @{
    3  = 2021
    4  = 58.14349653559799
    5  = 50.405857641278807
    6  = 65.995789757646122
    7  = 47.339416388110941
    8  = 44.18457369250482
    9  = 50.379263611270765
    10 = 54.819947539591084
    11 = 47.679920417302263
    12 = 61.861274529713718
    13 = 36.712395096811576
    14 = 26.872053459579295
    15 = 46.638444428499682
    16 = 51.155081745820631
    17 = 43.08338023862634
    18 = 58.934228062068456
    19 = 54.51979816984521
    20 = 52.474443936678909
    21 = 56.519551395440942
    22 = 46.970408642555192
    23 = 27.43769048802011
    24 = 66.104415920267911
    25 = 88.246666265390886
    26 = 71.914698721605745
    27 = 105.10059183863845
    28 = 63.980940123966526
    29 = 55.546587096180644
    30 = 73.505198287622903
    31 = 43.916363725083563
    32 = 40.980198843051781
    33 = 47.015458682814909
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 18).Value = $newValues[$row]
}

# 3) Move the active selection to T3, matching the author's saved view.
$ws.Range("T3").Select()
